$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (D1:M1): copy header formatting from A1, then set new text ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:M1").PasteSpecial(-4122) | Out-Null

$headers = @("particip","taxa_sucesso","arrecadado_sucesso","media_sucesso","std_sucesso","min_sucesso","max_sucesso","apoio_medio","contribuicoes","media_contribuicoes")
$col = 4
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# --- Data rows ---
# Columns: D=particip, E=taxa_sucesso, F=arrecadado_sucesso, G=media_sucesso,
#          H=std_sucesso, I=min_sucesso, J=max_sucesso, K=apoio_medio,
#          L=contribuicoes, M=media_contribuicoes
$data = @{
    2 = @(100, 15.15151515151515, 426.4616061876675, 85.29232123753351, 85.51030885495558, 7.154956142241136, 226.3900843036052, 18.5418089646812, 23, 4.6)
    3 = @(100, 9.929078014184398, 4736.523382339164, 338.3230987385117, 458.6831840070852, 27.62335886703489, 1809.09852121176, 16.91615493692558, 280, 20)
    4 = @(100, 14.66666666666667, 7327.98576013602, 666.180523648729, 1084.782975461776, 26.58043580770418, 3475.049171548047, 20.99709386858459, 349, 31.72727272727273)
    5 = @(100, 18.82352941176471, 3191.257392255826, 199.4535870159891, 234.097723686266, 10.31772032536115, 834.8528000913501, 17.06554755217019, 187, 11.6875)
    6 = @(100, 30.8411214953271, 15024.86326411449, 455.2988867913483, 1126.99584249591, 3.799754022893506, 5087.076865717208, 20.44199083552992, 735, 22.27272727272727)
    7 = @(100, 21.42857142857143, 2767.438420781004, 131.7827819419526, 138.1458157071481, 5.763382152582333, 538.4389998789497, 24.70927161411611, 112, 5.333333333333333)
    8 = @(100, 31.88405797101449, 4495.625203875862, 204.3466001761755, 345.0624243778461, 1.087396962410123, 1594.029696524064, 23.66118528355717, 190, 8.636363636363637)
    9 = @(100, 39.47368421052632, 5216.802725094768, 173.8934241698256, 195.9730930324503, 2.022084306600051, 657.0789958678034, 15.71326122016496, 332, 11.06666666666667)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $row = [int]$r

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $vals[0]
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $vals[1]
    $eCell.NumberFormat = "0.00%"

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $vals[2]
    $fCell.NumberFormat = "R$ #,##0.00"

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $vals[3]
    $gCell.NumberFormat = "R$ #,##0.00"

    $ws.Cells.Item($row, 8).Value = $vals[4]
    $ws.Cells.Item($row, 8).Style = "Normal"
    $ws.Cells.Item($row, 9).Value = $vals[5]
    $ws.Cells.Item($row, 9).Style = "Normal"
    $ws.Cells.Item($row, 10).Value = $vals[6]
    $ws.Cells.Item($row, 10).Style = "Normal"
    $ws.Cells.Item($row, 11).Value = $vals[7]
    $ws.Cells.Item($row, 11).Style = "Normal"
    $ws.Cells.Item($row, 12).Value = $vals[8]
    $ws.Cells.Item($row, 12).Style = "Normal"
    $ws.Cells.Item($row, 13).Value = $vals[9]
    $ws.Cells.Item($row, 13).Style = "Normal"
}